$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: min row ---
$ws.Range("D28").Value = "min"
$ws.Range("D28").VerticalAlignment = -4108

$ws.Range("E28").Formula = "=MIN(E3:E22)"
$ws.Range("F28:J28").Formula = "=MIN(F3:F22)"
$ws.Range("I28").ClearContents()

$row28 = $ws.Range("E28:J28")
$row28.NumberFormat = "0.000"
$row28.VerticalAlignment = -4108

# --- Row 29: max row ---
$ws.Range("D29").Value = "max"
$ws.Range("D29").VerticalAlignment = -4108

$ws.Range("E29").Formula = "=MAX(E3:E22)"
$ws.Range("F29:J29").Formula = "=MAX(F3:F22)"
$ws.Range("I29").ClearContents()

$row29 = $ws.Range("E29:J29")
$row29.NumberFormat = "0.000"
$row29.VerticalAlignment = -4108

# --- Selection / active cell as left by the author ---
$null = $ws.Range("E33").Select()
